{"js": "// Replace the division-problem text in the worksheet table.\n// The table has 20 rows x 5 columns, but only rows 0, 4, 8, 12, 16\n// (0-indexed) contain the actual \"NN\u00f7N=\" problems; the rest are blank\n// spacer rows. Each content row has exactly 5 cells (columns 0-4).\n// We update the text of each cell in-place (preserving paragraph /\n// run formatting) by replacing the paragraph's text.\n\nconst replacements = [\n  [0, 0, \"65\u00f77=\", \"56\u00f77=\"],\n  [0, 1, \"90\u00f74=\", \"34\u00f75=\"],\n  [0, 2, \"55\u00f77=\", \"90\u00f72=\"],\n  [0, 3, \"37\u00f74=\", \"32\u00f75=\"],\n  [0, 4, \"20\u00f78=\", \"56\u00f77=\"],\n  [4, 0, \"93\u00f76=\", \"74\u00f74=\"],\n  [4, 1, \"23\u00f79=\", \"29\u00f73=\"],\n  [4, 2, \"76\u00f77=\", \"94\u00f72=\"],\n  [4, 3, \"51\u00f74=\", \"78\u00f72=\"],\n  [4, 4, \"63\u00f77=\", \"61\u00f79=\"],\n  [8, 0, \"53\u00f77=\", \"19\u00f77=\"],\n  [8, 1, \"34\u00f76=\", \"27\u00f75=\"],\n  [8, 2, \"96\u00f79=\", \"13\u00f75=\"],\n  [8, 3, \"22\u00f78=\", \"86\u00f78=\"],\n  [8, 4, \"95\u00f78=\", \"94\u00f79=\"],\n  [12, 0, \"78\u00f77=\", \"51\u00f73=\"],\n  [12, 1, \"93\u00f76=\", \"94\u00f78=\"],\n  [12, 2, \"52\u00f75=\", \"87\u00f72=\"],\n  [12, 3, \"66\u00f72=\", \"21\u00f79=\"],\n  [12, 4, \"39\u00f72=\", \"21\u00f75=\"],\n  [16, 0, \"71\u00f72=\", \"89\u00f79=\"],\n  [16, 1, \"66\u00f72=\", \"36\u00f75=\"],\n  [16, 2, \"78\u00f77=\", \"40\u00f74=\"],\n  [16, 3, \"99\u00f74=\", \"20\u00f74=\"],\n  [16, 4, \"81\u00f74=\", \"54\u00f79=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Load every target cell's first paragraph text up front.\nconst cellParas = replacements.map(([row, col]) => {\n  const cell = table.getCell(row, col);\n  const para = cell.body.paragraphs.getFirst();\n  para.load(\"text\");\n  return para;\n});\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, , oldText, newText] = replacements[i];\n  const para = cellParas[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Unexpected cell text at index ${i}: expected \"${oldText}\" but found \"${para.text}\"`\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the division-problem text in the worksheet table.\n# The table has 20 rows x 5 columns, but only rows 1, 5, 9, 13, 17\n# (1-indexed, as COM table rows are) contain the actual \"NN\u00f7N=\" problems;\n# the rest are blank spacer rows. Each content row has exactly 5 cells\n# (columns 1-5). We update each cell's Range.Text in place, which\n# preserves the surrounding paragraph / run formatting.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @(1, 1, \"65\u00f77=\", \"56\u00f77=\"),\n    @(1, 2, \"90\u00f74=\", \"34\u00f75=\"),\n    @(1, 3, \"55\u00f77=\", \"90\u00f72=\"),\n    @(1, 4, \"37\u00f74=\", \"32\u00f75=\"),\n    @(1, 5, \"20\u00f78=\", \"56\u00f77=\"),\n    @(5, 1, \"93\u00f76=\", \"74\u00f74=\"),\n    @(5, 2, \"23\u00f79=\", \"29\u00f73=\"),\n    @(5, 3, \"76\u00f77=\", \"94\u00f72=\"),\n    @(5, 4, \"51\u00f74=\", \"78\u00f72=\"),\n    @(5, 5, \"63\u00f77=\", \"61\u00f79=\"),\n    @(9, 1, \"53\u00f77=\", \"19\u00f77=\"),\n    @(9, 2, \"34\u00f76=\", \"27\u00f75=\"),\n    @(9, 3, \"96\u00f79=\", \"13\u00f75=\"),\n    @(9, 4, \"22\u00f78=\", \"86\u00f78=\"),\n    @(9, 5, \"95\u00f78=\", \"94\u00f79=\"),\n    @(13, 1, \"78\u00f77=\", \"51\u00f73=\"),\n    @(13, 2, \"93\u00f76=\", \"94\u00f78=\"),\n    @(13, 3, \"52\u00f75=\", \"87\u00f72=\"),\n    @(13, 4, \"66\u00f72=\", \"21\u00f79=\"),\n    @(13, 5, \"39\u00f72=\", \"21\u00f75=\"),\n    @(17, 1, \"71\u00f72=\", \"89\u00f79=\"),\n    @(17, 2, \"66\u00f72=\", \"36\u00f75=\"),\n    @(17, 3, \"78\u00f77=\", \"40\u00f74=\"),\n    @(17, 4, \"99\u00f74=\", \"20\u00f74=\"),\n    @(17, 5, \"81\u00f74=\", \"54\u00f79=\")\n)\n\nforeach ($r in $replacements) {\n    $row = $r[0]\n    $col = $r[1]\n    $oldText = $r[2]\n    $newText = $r[3]\n\n    $cell = $t.Cell($row, $col)\n    $cellRange = $cell.Range\n    $cellRange.MoveEnd(1, -1) | Out-Null  # trim trailing end-of-cell mark\n    if ($cellRange.Text -ne $oldText) {\n        throw \"Unexpected cell text at row $row col ${col}: expected '$oldText' but found '$($cellRange.Text)'\"\n    }\n    $cellRange.Text = $newText\n}\n"}
